$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column A (xF) values from row 2 to row 17 to 0.6
$ws.Range("A2:A17").Value = 0.6

# Move selection to D23 as reflected in the saved view state
$ws.Range("D23").Select()
